$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scores")

# --- Mirror header row (A1:D1) into F1:I1, preserving style ---
$ws.Range("A1:D1").Copy($ws.Range("F1"))

# --- Mirror model-name column (A2:A6) into F2:F6 ---
$ws.Range("A2:A6").Copy($ws.Range("F2"))

# --- ROUND() comparison columns G:I ---
$ws.Range("G2").Formula = "=ROUND(B2,2)"
$ws.Range("H2:I2").Formula = "=ROUND(C2,2)"
$ws.Range("G3:G6").Formula = "=ROUND(B3,2)"
$ws.Range("H3:H6").Formula = "=ROUND(C3,2)"
$ws.Range("I3:I6").Formula = "=ROUND(D3,2)"

# --- Column widths for the new comparison table ---
$ws.Columns.Item(7).ColumnWidth = 14.1666667
$ws.Columns.Item(9).ColumnWidth = 16.1666667

# --- Page setup for the Scores sheet ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection on the Scores sheet moves to the new table ---
[void]$ws.Range("F1:I6").Select()

# --- "dimension reduction" sheet becomes the active tab ---
$ws2 = $wb.Worksheets.Item("dimension reduction")
[void]$ws2.Activate()
